# Updates cryptos list values (prices / 1h volume %) and fixes the
# EnergySwap / Monero row ordering, per the Sat Feb 17 20:58:15 UTC 2024
# GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would NOT be re-parsed as a plain number by Excel
# (percentages, multi-dot "thousands" style prices, the coin name/link
# swap for rows 43-44, etc.) -- a straight .Value assignment is safe here
# and keeps the cell's existing (default) style untouched.
$plainUpdates = @{
    'D2' = '51.509.61'
    'E2' = '  -0.62%  '
    'D3' = '2.778.48'
    'E3' = '  -0.05%  '
    'E4' = '  -0.07%  '
    'E5' = '  -1.33%  '
    'E6' = '  -1.23%  '
    'E7' = '  -3.10%  '
    'E8' = '  +0.05%  '
    'E9' = '  +0.44%  '
    'E10' = '  -0.84%  '
    'E11' = '  +3.26%  '
    'E12' = '  -1.62%  '
    'E13' = '  +3.35%  '
    'E14' = '  +0.84%  '
    'D15' = '3.212.65'
    'D16' = '2.783.49'
    'E16' = '  +0.07%  '
    'E17' = '  -1.09%  '
    'D18' = '51.461.01'
    'E18' = '  -0.59%  '
    'E19' = '  +3.47%  '
    'E20' = '  -1.39%  '
    'E21' = '  +0.61%  '
    'D22' = '0.0₃0962'
    'E22' = '  -1.47%  '
    'E23' = '  -0.25%  '
    'E24' = '  -2.75%  '
    'E25' = '  -0.73%  '
    'E26' = '  -0.09%  '
    'E27' = '  -2.03%  '
    'E29' = '  +0.70%  '
    'E30' = '  -0.04%  '
    'E31' = '  +7.24%  '
    'E32' = '  +9.69%  '
    'E33' = '  +0.83%  '
    'E34' = '  -2.37%  '
    'E35' = '  +6.17%  '
    'E36' = '  -2.00%  '
    'E37' = '  -0.13%  '
    'E38' = '  +2.03%  '
    'E39' = '  -2.99%  '
    'E40' = '  -1.82%  '
    'E41' = '  -0.16%  '
    'E42' = '  -0.92%  '
    'B43' = 'Monero'
    'C43' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E43' = '  -1.38%  '
    'B44' = 'EnergySwap'
    'C44' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E44' = '  +0.28%  '
    'D46' = '2.109.28'
    'E46' = '  +2.08%  '
    'E47' = '  +0.87%  '
    'E48' = '  +6.87%  '
    'E49' = '  -4.75%  '
    'E50' = '  -3.20%  '
    'E51' = '  +7.80%  '
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Price cells in column D whose new text DOES look like a plain number
# (e.g. "108.40"). A bare .Value assignment would get auto-coerced to a
# Double and silently drop the trailing zero / formatting, so each one is
# forced to Text just long enough to hold the literal string, then the
# cell style is reset back to Normal/General to match the original
# (unstyled) cells.
$numericTextUpdates = @{
    'D5' = '352.71'
    'D6' = '108.40'
    'D7' = '0.548'
    'D9' = '0.597'
    'D10' = '39.71'
    'D12' = '0.0836'
    'D13' = '20.10'
    'D14' = '7.66'
    'D17' = '0.924'
    'D19' = '7.68'
    'D21' = '13.14'
    'D23' = '69.90'
    'D24' = '265.97'
    'D27' = '26.10'
    'D29' = '10.23'
    'D30' = '2.21'
    'D31' = '36.33'
    'D32' = '6.25'
    'D33' = '51.94'
    'D34' = '0.0453'
    'D35' = '5.57'
    'D36' = '0.0827'
    'D37' = '0.999'
    'D38' = '18.44'
    'D41' = '2.52'
    'D43' = '119.94'
    'D44' = '22.03'
    'D47' = '3.27'
    'D48' = '2.32'
    'D49' = '5.42'
    'D50' = '0.904'
}

foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
    $cell.Style = "Normal"
}
